$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update HUMAN_E (column G) scores from 0 to 1 for re-evaluated rows
$ws.Cells.Item(287, 7).Value = 1
$ws.Cells.Item(288, 7).Value = 1
$ws.Cells.Item(289, 7).Value = 1
$ws.Cells.Item(292, 7).Value = 1
$ws.Cells.Item(293, 7).Value = 1
$ws.Cells.Item(294, 7).Value = 1
$ws.Cells.Item(296, 7).Value = 1
$ws.Cells.Item(298, 7).Value = 1
$ws.Cells.Item(300, 7).Value = 1
$ws.Cells.Item(301, 7).Value = 1
$ws.Cells.Item(302, 7).Value = 1
$ws.Cells.Item(307, 7).Value = 1
$ws.Cells.Item(308, 7).Value = 1
$ws.Cells.Item(309, 7).Value = 1
$ws.Cells.Item(311, 7).Value = 1
$ws.Cells.Item(313, 7).Value = 1
$ws.Cells.Item(314, 7).Value = 1
$ws.Cells.Item(315, 7).Value = 1
$ws.Cells.Item(321, 7).Value = 1
$ws.Cells.Item(322, 7).Value = 1
$ws.Cells.Item(323, 7).Value = 1
$ws.Cells.Item(324, 7).Value = 1
$ws.Cells.Item(325, 7).Value = 1

# Add human-evaluator comments (column H) explaining the score/notes
$ws.Cells.Item(294, 8).Value = "Qualche errore sintattico, però va bene"
$ws.Cells.Item(295, 8).Value = "Non fa l'ultima operazione"
$ws.Cells.Item(300, 8).Value = "Qualche errore sintattico, però va bene"
$ws.Cells.Item(303, 8).Value = "Non genera il clock come richiesto"
$ws.Cells.Item(304, 8).Value = "Completamente sbagliato"
$ws.Cells.Item(306, 8).Value = "non digliara il segnale ytemp a 4 bit"
$ws.Cells.Item(314, 8).Value = "Inefficiente ma va bene"
$ws.Cells.Item(316, 8).Value = "Non dichiara la libreria standard alla fine"
$ws.Cells.Item(318, 8).Value = "Non equivalenti per quello che è stato richiesto"
$ws.Cells.Item(319, 8).Value = "Noin fa la concatenazione "
$ws.Cells.Item(326, 8).Value = "Non fa l'inizializzazione come richiesto"

# Update the sheet view: scroll position and active selection
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 297
$win.ScrollColumn = 1
$ws.Range("H325").Select()
